$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1) Swap the match data between rows 89 and 90 (keep A,B,C,D,E,G,M,Q,U as-is)
# ---------------------------------------------------------------
$cols = @(6, 8, 9, 10, 11, 12, 14, 15, 16, 18, 19, 20, 22)  # F,H,I,J,K,L,N,O,P,R,S,T,V

foreach ($c in $cols) {
    $v89 = $ws.Cells.Item(89, $c).Value2
    $v90 = $ws.Cells.Item(90, $c).Value2
    $ws.Cells.Item(89, $c).Value2 = $v90
    $ws.Cells.Item(90, $c).Value2 = $v89
}

# ---------------------------------------------------------------
# 2) Append three new match rows (102, 103, 104)
# ---------------------------------------------------------------
function Set-Row {
    param($r, $idx, $home, $homeGoals, $away, $awayGoals, $dateSerial,
          $homeOpen, $homeOpenDt, $homeClose, $homeCloseDt,
          $drawOpen, $drawOpenDt, $drawClose, $drawCloseDt,
          $awayOpen, $awayOpenDt, $awayClose, $awayCloseDt, $url)

    $ws.Cells.Item($r, 1).Value2 = $idx
    $ws.Cells.Item($r, 2).Value2 = "turkey"
    $ws.Cells.Item($r, 3).Value2 = "1-lig"
    $ws.Cells.Item($r, 4).Value2 = "2023-2024"
    $ws.Cells.Item($r, 5).Value2 = $dateSerial
    $ws.Cells.Item($r, 6).Value2 = $home
    $ws.Cells.Item($r, 7).Value2 = $homeGoals
    $ws.Cells.Item($r, 8).Value2 = $away
    $ws.Cells.Item($r, 9).Value2 = $awayGoals
    $ws.Cells.Item($r, 10).Value2 = $homeOpen
    $ws.Cells.Item($r, 11).Value2 = $homeOpenDt
    $ws.Cells.Item($r, 12).Value2 = $homeClose
    $ws.Cells.Item($r, 13).Value2 = $homeCloseDt
    $ws.Cells.Item($r, 14).Value2 = $drawOpen
    $ws.Cells.Item($r, 15).Value2 = $drawOpenDt
    $ws.Cells.Item($r, 16).Value2 = $drawClose
    $ws.Cells.Item($r, 17).Value2 = $drawCloseDt
    $ws.Cells.Item($r, 18).Value2 = $awayOpen
    $ws.Cells.Item($r, 19).Value2 = $awayOpenDt
    $ws.Cells.Item($r, 20).Value2 = $awayClose
    $ws.Cells.Item($r, 21).Value2 = $awayCloseDt
    $ws.Cells.Item($r, 22).Value2 = $url

    $ws.Cells.Item(101, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item(101, 5).Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)
}

Set-Row 102 101 "Keciorengucu" 1 "Bodrumspor" 1 45241.47916666666 `
    2.97 "05/11/2023 11:42" 3.3 "11/11/2023 11:22" `
    3.07 "05/11/2023 11:42" 3.43 "11/11/2023 11:22" `
    2.49 "05/11/2023 11:42" 2.22 "11/11/2023 11:22" `
    "https://www.betexplorer.com/football/turkey/1-lig/keciorengucu-bodrumspor/W0Q7oRXj/"

Set-Row 103 102 "Tuzlaspor" 3 "Giresunspor" 1 45241.47916666666 `
    1.84 "06/11/2023 18:12" 1.97 "11/11/2023 11:27" `
    3.6 "06/11/2023 18:12" 3.44 "11/11/2023 11:27" `
    4.27 "06/11/2023 18:12" 4.02 "11/11/2023 11:27" `
    "https://www.betexplorer.com/football/turkey/1-lig/tuzlaspor-giresunspor/OfDg9UPc/"

Set-Row 104 103 "Sanliurfaspor" 0 "Umraniyespor" 1 45241.58333333334 `
    2.48 "06/11/2023 18:12" 3.09 "11/11/2023 13:59" `
    3.24 "06/11/2023 18:12" 3.56 "11/11/2023 13:59" `
    2.92 "06/11/2023 18:12" 2.27 "11/11/2023 13:59" `
    "https://www.betexplorer.com/football/turkey/1-lig/sanliurfaspor-umraniyespor/CGqyimAT/"
